$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 175: fill in the remaining figures for 09-09-2021 ---
$ws.Range("B175").Value = 180000
$ws.Range("C175").Value = 1.5
$ws.Range("D175").Value = 1.5
$ws.Range("E175").Value = 1.5
$ws.Range("F175").Value = 3
$ws.Range("G175").Value = 1.5

# --- Row 176: 10-09-2021 ---
# NOTE: typing "10-09-2021" directly makes Excel's smart-entry parse it as a
# date (10/9/2021), which would store it as a serial number with a new
# number-format style instead of the plain text label used throughout column
# A. Build it as a text formula in a scratch cell and paste-special just the
# values in, which keeps it a plain shared-string text cell with no special
# formatting (matching how the rest of the date labels are stored).
$ws.Range("ZZ1").Formula = "=""10-09-2021"""
$ws.Range("ZZ1").Copy()
$ws.Range("A176").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Range("B176").Value = 301000
$ws.Range("C176").Value = 1.5
$ws.Range("D176").Value = 1.5
$ws.Range("E176").Value = 1.5
$ws.Range("F176").Value = 7
$ws.Range("G176").Value = 1.5

# --- Row 177: 13-09-2021 ---
$ws.Range("A177").Value = "13-09-2021"
$ws.Range("B177").Value = 335000
$ws.Range("C177").Value = 1.5
$ws.Range("D177").Value = 1.5
$ws.Range("E177").Value = 1.5
$ws.Range("F177").Value = 6
$ws.Range("G177").Value = 1.5

# --- Row 178: 14-09-2021 ---
$ws.Range("A178").Value = "14-09-2021"
$ws.Range("B178").Value = 445000
$ws.Range("C178").Value = 1.5
$ws.Range("D178").Value = 1.5
$ws.Range("E178").Value = 1.5
$ws.Range("F178").Value = 7
$ws.Range("G178").Value = 1.5

# --- Row 179: 15-09-2021 ---
$ws.Range("A179").Value = "15-09-2021"
$ws.Range("B179").Value = 200000
$ws.Range("C179").Value = 1.5
$ws.Range("D179").Value = 1.5
$ws.Range("E179").Value = 1.5
$ws.Range("F179").Value = 4
$ws.Range("G179").Value = 1.5

# --- Row 180: 16-09-2021 ---
$ws.Range("A180").Value = "16-09-2021"
$ws.Range("B180").Value = 155000
$ws.Range("C180").Value = 1.5
$ws.Range("D180").Value = 1.5
$ws.Range("E180").Value = 1.5
$ws.Range("F180").Value = 6
$ws.Range("G180").Value = 1.5

# --- Row 181: 20-09-2021 ---
$ws.Range("A181").Value = "20-09-2021"
$ws.Range("G181").Value = 1.5
